# "More 8 Start 9" - corrections to Versuch 6 Berechnungen (Protokoll 8)
#
# - B17 / B20: offset constant 0.14994 -> 0.04994
# - C43: measured value correction 2.03495 -> 2.033
# - M43:M47: measurement-uncertainty formula now weights K/E and L/F terms
#   by a factor of 3 (SQRT((I/C)^2+(J/D)^2+(3*K/E)^2+(3*L/F)^2))
# - G45: shared-formula anchor range narrows to G45:G46
# - view repositioned to the bottom of the sheet, selection on M43:M47

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- corrected offset constants -----------------------------------------
$ws.Range("B17").Formula = "=B16+0.04994"
$ws.Range("B20").Formula = "=B19+0.04994"

# --- corrected measured diameter ----------------------------------------
$ws.Range("C43").Value = 2.033

# --- corrected uncertainty-propagation formula (rows 43-47) -------------
$ws.Range("M43").Formula = "=SQRT((I43/C43)^2+(J43/D43)^2+(3*K43/E43)^2+(3*L43/F43)^2)"
$ws.Range("M44:M47").Formula = "=SQRT((I44/C44)^2+(J44/D44)^2+(3*K44/E44)^2+(3*L44/F44)^2)"

# --- view / selection -----------------------------------------------------
$ws.Range("M43:M47").Select()
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 6
